# Update odds values in "Jogos da Semana" workbook (FlashScore export).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5
$ws.Range("Q5").Value = 2.05
$ws.Range("R5").Value = 1.8

# Row 14
$ws.Range("M14").Value = 1.01
$ws.Range("N14").Value = 7.6

# Row 17
$ws.Range("G17").Value = 1.35
$ws.Range("H17").Value = 5.2
$ws.Range("I17").Value = 5.7
$ws.Range("J17").Value = 1.7
$ws.Range("K17").Value = 2.95
$ws.Range("L17").Value = 5
$ws.Range("Q17").Value = 1.18
$ws.Range("R17").Value = 3.74
$ws.Range("U17").Value = 1.4
$ws.Range("V17").Value = 2.85
$ws.Range("X17").Value = 10.25
$ws.Range("Y17").Value = 8.5
$ws.Range("Z17").Value = 10.5
$ws.Range("AA17").Value = 8.75
$ws.Range("AB17").Value = 13
$ws.Range("AD17").Value = 11.5
$ws.Range("AE17").Value = 12.5
$ws.Range("AF17").Value = 28
$ws.Range("AG17").Value = 100
$ws.Range("AH17").Value = 29
$ws.Range("AI17").Value = 45
$ws.Range("AJ17").Value = 17
$ws.Range("AK17").Value = 90
$ws.Range("AL17").Value = 37
$ws.Range("AM17").Value = 27
$ws.Range("AN17").Value = 4.15
$ws.Range("AO17").Value = 6.1
$ws.Range("AQ17").Value = 13
$ws.Range("AR17").Value = 22
$ws.Range("AT17").Value = 5.3
$ws.Range("AU17").Value = 7
$ws.Range("AV17").Value = 32
$ws.Range("AW17").Value = 300
$ws.Range("AX17").Value = 8.75
$ws.Range("AY17").Value = 27
$ws.Range("AZ17").Value = 22
$ws.Range("BA17").Value = 120
$ws.Range("BB17").Value = 100
$ws.Range("BC17").Value = 150
